# Add a new "2023" column (L) to the 16.6.2.1a worksheet, mirroring the
# existing "2022" column (K): same per-row formatting, new values, taller
# data rows, and a refreshed used range / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Write the new 2023 values into column L (row 4 = year header).
# ---------------------------------------------------------------------
$values = @{
    4  = 2023
    5  = 33.991563806511245
    6  = 43.352267904134116
    7  = 46.016552065013244
    8  = 57.950845675564537
    9  = 46.481788079470263
    10 = 45.080578284701389
    11 = 39.506289942950417
    12 = 26.964612178240138
    13 = 15.46142526802614
    14 = 33.453947368420813
}

foreach ($row in $values.Keys) {
    $srcCell = $ws.Cells.Item($row, 11)   # K<row>  (2022 column, same formatting we want to mirror)
    $dstCell = $ws.Cells.Item($row, 12)   # L<row>  (new 2023 column)

    # Copy formatting only (keeps the destination's own value slot free so
    # we can set the real 2023 number right after).
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $dstCell.Value = $values[$row]
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Taller rows for the data table (4 through 14) to fit the extra
#    column / updated layout; header row 1 also grows.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 43.5

$dataRows = 4..14
foreach ($r in $dataRows) {
    $ws.Rows.Item($r).RowHeight = 14.25
}

# ---------------------------------------------------------------------
# 3. Reset the saved selection back to the top-left cell.
# ---------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
